$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $savedStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $savedStyle
}

$ws.Range("D2").Value = '26.651.80'
$ws.Range("E2").Value = '  -0.43%  '

$ws.Range("D3").Value = '1.596.47'
$ws.Range("E3").Value = '  -0.25%  '

$ws.Range("E4").Value = '  -0.01%  '

Set-TextValue $ws 'D5' '211.35'
$ws.Range("E5").Value = '  +0.56%  '

Set-TextValue $ws 'D6' '0.510'
$ws.Range("E6").Value = '  +1.17%  '

$ws.Range("E7").Value = '  +0.03%  '

$ws.Range("E8").Value = '  +0.05%  '

$ws.Range("E9").Value = '  -0.83%  '

Set-TextValue $ws 'D10' '19.67'
$ws.Range("E10").Value = '  +0.30%  '

$ws.Range("E11").Value = '  +0.06%  '

$ws.Range("D12").Value = '1.821.25'
$ws.Range("E12").Value = '  -0.61%  '

$ws.Range("D13").Value = '1.573.51'
$ws.Range("E13").Value = '  -2.40%  '

$ws.Range("E14").Value = '  -0.45%  '

$ws.Range("E15").Value = '  -1.41%  '

Set-TextValue $ws 'D16' '64.92'
$ws.Range("E16").Value = '  +2.33%  '

$ws.Range("D17").Value = '26.644.90'
$ws.Range("E17").Value = '  -0.44%  '

$ws.Range("E18").Value = '  +0.19%  '

Set-TextValue $ws 'D19' '209.48'
$ws.Range("E19").Value = '  +0.13%  '

$ws.Range("E20").Value = '  +0.10%  '

Set-TextValue $ws 'D21' '6.75'
$ws.Range("E21").Value = '  +0.41%  '

$ws.Range("E22").Value = '  -0.01%  '

$ws.Range("E23").Value = '  -1.16%  '

$ws.Range("E24").Value = '  +0.81%  '

Set-TextValue $ws 'D25' '146.37'
$ws.Range("E25").Value = '  +0.13%  '

$ws.Range("E26").Value = '  -0.02%  '

Set-TextValue $ws 'D27' '7.16'
$ws.Range("E27").Value = '  -4.08%  '

$ws.Range("E28").Value = '  +2.44%  '

$ws.Range("E29").Value = '  -0.03%  '

$ws.Range("E30").Value = '  +0.74%  '

$ws.Range("E31").Value = '  +0.18%  '

Set-TextValue $ws 'D32' '3.23'
$ws.Range("E32").Value = '  -0.50%  '

$ws.Range("E33").Value = '  -0.44%  '

$ws.Range("E34").Value = '  -0.68%  '

$ws.Range("D35").Value = '1.298.55'
$ws.Range("E35").Value = '  -1.06%  '

$ws.Range("E36").Value = '  +0.51%  '

$ws.Range("E37").Value = '  -2.24%  '

$ws.Range("E38").Value = '  -0.69%  '

Set-TextValue $ws 'D39' '0.840'
$ws.Range("E39").Value = '  +2.65%  '

$ws.Range("E40").Value = '  +0.05%  '

$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue $ws 'D41' '5.37'
$ws.Range("E41").Value = '  +1.84%  '

$ws.Range("B42").Value = 'MXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws 'D42' '2.20'
$ws.Range("E42").Value = '  +0.88%  '

$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws 'D43' '0.788'
$ws.Range("E43").Value = '  +0.00%  '

Set-TextValue $ws 'D44' '63.83'
$ws.Range("E44").Value = '  +1.82%  '

$ws.Range("D45").Value = '1.733.51'
$ws.Range("E45").Value = '  -0.64%  '

$ws.Range("E46").Value = '  +8.98%  '

Set-TextValue $ws 'D47' '90.08'
$ws.Range("E47").Value = '  +1.25%  '

Set-TextValue $ws 'D48' '1.62'
$ws.Range("E48").Value = '  +0.78%  '

Set-TextValue $ws 'D49' '0.0998'
$ws.Range("E49").Value = '  +2.39%  '

$ws.Range("E50").Value = '  -1.16%  '

$ws.Range("E51").Value = '  +0.35%  '
